# Swap the match-data (everything except id/Div/Date columns A, C, D) between
# pairs of rows that were mis-ordered in the source feed. Each pair below
# represents the same two underlying matches whose rows need their
# B (id) and E:AD (HomeTeam..PL_AhUnder) values exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(35, 36),
    @(64, 65),
    @(73, 74),
    @(116, 117),
    @(118, 121),
    @(156, 157),
    @(158, 159),
    @(194, 195),
    @(213, 214)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1" + ":AD$r1")
    $range2 = $ws.Range("B$r2" + ":AD$r2")

    $v1 = $range1.Value2
    $v2 = $range2.Value2

    $range1.Value = $v2
    $range2.Value = $v1
}
